$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("N8").Value = "xtndthre jxtym"

$ws.Range("I13").Select()
